$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44893
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 900
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 13444
$ws.Range("N2").Value = '$/caja 13 kilos'
$ws.Range("P2").Value = 1034
$ws.Range("Q2").Value = 13

$ws.Range("D3").Value = 44984
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("N3").Value = '$/caja 13 kilos'
$ws.Range("P3").Value = 1269
$ws.Range("Q3").Value = 13

$ws.Range("D4").Value = 45245
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("N4").Value = '$/caja 13 kilos'
$ws.Range("P4").Value = 1346
$ws.Range("Q4").Value = 13

$ws.Range("D5").Value = 45141
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = 16550
$ws.Range("N5").Value = '$/caja 13 kilos'
$ws.Range("P5").Value = 1273
$ws.Range("Q5").Value = 13

$ws.Range("D6").Value = 45154
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("N6").Value = '$/caja 13 kilos'
$ws.Range("P6").Value = 1346
$ws.Range("Q6").Value = 13

$ws.Range("D7").Value = 45155
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 16500
$ws.Range("N7").Value = '$/caja 13 kilos'
$ws.Range("P7").Value = 1269
$ws.Range("Q7").Value = 13

$ws.Range("D8").Value = 44469
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13500
$ws.Range("N8").Value = '$/caja 13 kilos'
$ws.Range("P8").Value = 1038
$ws.Range("Q8").Value = 13

$ws.Range("D9").Value = 44406
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("N9").Value = '$/caja 13 kilos'
$ws.Range("P9").Value = 1346
$ws.Range("Q9").Value = 13

$ws.Range("D10").Value = 44592
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("N10").Value = '$/caja 13 kilos'
$ws.Range("P10").Value = 962
$ws.Range("Q10").Value = 13

$ws.Range("D11").Value = 44320
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 19000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19500
$ws.Range("N11").Value = '$/caja 13 kilos'
$ws.Range("P11").Value = 1500
$ws.Range("Q11").Value = 13

$ws.Range("D12").Value = 45028
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("N12").Value = '$/caja 13 kilos'
$ws.Range("P12").Value = 1115
$ws.Range("Q12").Value = 13

$ws.Range("D13").Value = 45049
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 13500
$ws.Range("N13").Value = '$/caja 13 kilos'
$ws.Range("P13").Value = 1038
$ws.Range("Q13").Value = 13

$ws.Range("D14").Value = 45238
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 16000
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 16500
$ws.Range("N14").Value = '$/caja 13 kilos'
$ws.Range("P14").Value = 1269
$ws.Range("Q14").Value = 13

$ws.Range("D15").Value = 45096
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 750
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14600
$ws.Range("N15").Value = '$/caja 13 kilos'
$ws.Range("P15").Value = 1123
$ws.Range("Q15").Value = 13

$ws.Range("D16").Value = 44972
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 350
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17429
$ws.Range("N16").Value = '$/caja 15 kilos'
$ws.Range("P16").Value = 1162
$ws.Range("Q16").Value = 15

$ws.Range("D17").Value = 45092
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 13000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 13500
$ws.Range("N17").Value = '$/caja 13 kilos'
$ws.Range("P17").Value = 1038
$ws.Range("Q17").Value = 13

$ws.Range("D18").Value = 44616
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19500
$ws.Range("N18").Value = '$/caja 13 kilos'
$ws.Range("P18").Value = 1500
$ws.Range("Q18").Value = 13

$ws.Range("D19").Value = 44914
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14400
$ws.Range("N19").Value = '$/caja 13 kilos'
$ws.Range("P19").Value = 1108
$ws.Range("Q19").Value = 13

$ws.Range("D20").Value = 44890
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("N20").Value = '$/caja 13 kilos'
$ws.Range("P20").Value = 1115
$ws.Range("Q20").Value = 13

$ws.Range("D21").Value = 44910
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("N21").Value = '$/caja 13 kilos'
$ws.Range("P21").Value = 1115
$ws.Range("Q21").Value = 13

$ws.Range("D22").Value = 44988
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17400
$ws.Range("N22").Value = '$/caja 13 kilos'
$ws.Range("P22").Value = 1338
$ws.Range("Q22").Value = 13

$ws.Range("D23").Value = 44918
$ws.Range("I23").Value = 'Segunda'
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12750
$ws.Range("N23").Value = '$/caja 13 kilos'
$ws.Range("P23").Value = 981
$ws.Range("Q23").Value = 13

$ws.Range("D24").Value = 44397
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 140
$ws.Range("K24").Value = 12500
$ws.Range("L24").Value = 13000
$ws.Range("M24").Value = 12750
$ws.Range("N24").Value = '$/caja 13 kilos'
$ws.Range("P24").Value = 981
$ws.Range("Q24").Value = 13

$ws.Range("D25").Value = 45100
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 15500
$ws.Range("N25").Value = '$/caja 13 kilos'
$ws.Range("P25").Value = 1192
$ws.Range("Q25").Value = 13

$ws.Range("D26").Value = 44943
$ws.Range("I26").Value = 'Segunda'
$ws.Range("J26").Value = 350
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14429
$ws.Range("N26").Value = '$/caja 13 kilos'
$ws.Range("P26").Value = 1110
$ws.Range("Q26").Value = 13

$ws.Range("D27").Value = 45140
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 16500
$ws.Range("N27").Value = '$/caja 13 kilos'
$ws.Range("P27").Value = 1269
$ws.Range("Q27").Value = 13

$ws.Range("D28").Value = 44159
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 23000
$ws.Range("L28").Value = 24000
$ws.Range("M28").Value = 23500
$ws.Range("N28").Value = '$/caja 13 kilos'
$ws.Range("P28").Value = 1808
$ws.Range("Q28").Value = 13

$ws.Range("D29").Value = 44580
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 11000
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = 11500
$ws.Range("N29").Value = '$/caja 13 kilos'
$ws.Range("P29").Value = 885
$ws.Range("Q29").Value = 13

$ws.Range("D30").Value = 44832
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 13000
$ws.Range("L30").Value = 14000
$ws.Range("M30").Value = 13500
$ws.Range("N30").Value = '$/caja 13 kilos'
$ws.Range("P30").Value = 1038
$ws.Range("Q30").Value = 13

$ws.Range("D31").Value = 44389
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 120
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12500
$ws.Range("N31").Value = '$/caja 13 kilos'
$ws.Range("P31").Value = 962
$ws.Range("Q31").Value = 13

$ws.Range("D32").Value = 44855
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("N32").Value = '$/caja 13 kilos'
$ws.Range("P32").Value = 769
$ws.Range("Q32").Value = 13

$ws.Range("D33").Value = 44764
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 12500
$ws.Range("N33").Value = '$/caja 13 kilos'
$ws.Range("P33").Value = 962
$ws.Range("Q33").Value = 13

$ws.Range("D34").Value = 44379
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 120
$ws.Range("K34").Value = 12000
$ws.Range("L34").Value = 13000
$ws.Range("M34").Value = 12667
$ws.Range("N34").Value = '$/caja 13 kilos'
$ws.Range("P34").Value = 974
$ws.Range("Q34").Value = 13

$ws.Range("D35").Value = 45239
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 800
$ws.Range("K35").Value = 16000
$ws.Range("L35").Value = 17000
$ws.Range("M35").Value = 16500
$ws.Range("N35").Value = '$/caja 13 kilos'
$ws.Range("P35").Value = 1269
$ws.Range("Q35").Value = 13
